$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "43.124.23"
$ws.Range("E2").Value = "  +0.26%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "2.315.86"
$ws.Range("E3").Value = "  +0.20%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.02%  "

# Row 5 - BNB
$ws.Range("D5").Value = "302.17"
$ws.Range("E5").Value = "  -0.38%  "

# Row 6 - Solana
$ws.Range("D6").Value = "99.08"
$ws.Range("E6").Value = "  -1.87%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  +1.11%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  -0.01%  "

# Row 9 - Cardano
$ws.Range("E9").Value = "  +1.72%  "

# Row 10 - Avalanche
$ws.Range("D10").Value = "35.85"
$ws.Range("E10").Value = "  +1.66%  "

# Row 11 - Dogecoin
$ws.Range("E11").Value = "  -0.67%  "

# Row 12 - TRON
$ws.Range("E12").Value = "  -1.03%  "

# Row 13 - Chainlink
$ws.Range("E13").Value = "  +0.98%  "

# Row 14 - Polkadot
$ws.Range("D14").Value = "6.92"
$ws.Range("E14").Value = "  -0.09%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").Value = "2.675.60"
$ws.Range("E15").Value = "  -0.44%  "

# Row 16 - WrappedEther
$ws.Range("D16").Value = "2.284.88"
$ws.Range("E16").Value = "  -2.73%  "

# Row 17 - Polygon
$ws.Range("E17").Value = "  -2.67%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "43.029.73"
$ws.Range("E18").Value = "  +0.12%  "

# Row 19
$ws.Range("D19").Value = "13.56"
$ws.Range("E19").Value = "  +7.01%  "

# Row 20
$ws.Range("E20").Value = "  +0.72%  "

# Row 21
$ws.Range("D21").Value = "6.18"
$ws.Range("E21").Value = "  +0.16%  "

# Row 22
$ws.Range("D22").Value = "68.00"
$ws.Range("E22").Value = "  +0.17%  "

# Row 23
$ws.Range("D23").Value = "240.53"
$ws.Range("E23").Value = "  +1.33%  "

# Row 24
$ws.Range("E24").Value = "  -1.96%  "

# Row 25 - becomes PancakeSwap (was Dai)
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  -0.97%  "

# Row 26 - becomes Dai (was PancakeSwap)
$ws.Range("B26").Value = "Dai"
$ws.Range("C26").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.03%  "

# Row 27
$ws.Range("D27").Value = "25.07"

# Row 28
$ws.Range("D28").Value = "168.82"
$ws.Range("E28").Value = "  +0.61%  "

# Row 29
$ws.Range("D29").Value = "9.21"
$ws.Range("E29").Value = "  -0.33%  "

# Row 30
$ws.Range("E30").Value = "  -10.45%  "

# Row 31
$ws.Range("D31").Value = "33.55"
$ws.Range("E31").Value = "  -1.50%  "

# Row 32
$ws.Range("D32").Value = "5.24"
$ws.Range("E32").Value = "  +4.53%  "

# Row 33
$ws.Range("D33").Value = "4.90"
$ws.Range("E33").Value = "  +6.34%  "

# Row 34
$ws.Range("E34").Value = "  -0.12%  "

# Row 35
$ws.Range("D35").Value = "18.37"
$ws.Range("E35").Value = "  +8.02%  "

# Row 36
$ws.Range("E36").Value = "  -0.06%  "

# Row 37
$ws.Range("E37").Value = "  +0.26%  "

# Row 38
$ws.Range("E38").Value = "  +0.15%  "

# Row 39
$ws.Range("E39").Value = "  +0.88%  "

# Row 40
$ws.Range("D40").Value = "2.77"
$ws.Range("E40").Value = "  -2.01%  "

# Row 41
$ws.Range("E41").Value = "  -0.18%  "

# Row 42
$ws.Range("D42").Value = "1.999.27"
$ws.Range("E42").Value = "  -0.19%  "

# Row 43
$ws.Range("E43").Value = "  +0.45%  "

# Row 44
$ws.Range("D44").Value = "2.16"
$ws.Range("E44").Value = "  -5.91%  "

# Row 45
$ws.Range("D45").Value = "10.11"
$ws.Range("E45").Value = "  -1.18%  "

# Row 46
$ws.Range("D46").Value = "17.48"
$ws.Range("E46").Value = "  -0.04%  "

# Row 47
$ws.Range("E47").Value = "  +0.01%  "

# Row 48
$ws.Range("D48").Value = "54.91"
$ws.Range("E48").Value = "  -0.93%  "

# Row 49
$ws.Range("D49").Value = "75.25"
$ws.Range("E49").Value = "  +7.34%  "

# Row 50
$ws.Range("D50").Value = "2.540.32"
$ws.Range("E50").Value = "  +0.37%  "

# Row 51
$ws.Range("E51").Value = "  +0.39%  "
